$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# Reference cell with the default (unstyled) body format, used to reset
# style on cells we briefly mark as Text (see below) so no stray style
# index is left attached to the edited cell.
$plainStyle = $ws.Range("C2").Style

$ws.Range("D2").Value = "37.192.01"
$ws.Range("E2").Value = "  +1.70%  "
$ws.Range("D3").Value = "2.002.93"
$ws.Range("E3").Value = "  +2.48%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.43"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = "  +0.98%  "
$ws.Range("E6").Value = "  +2.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.34"
$ws.Range("D7").Style = $plainStyle
$ws.Range("E7").Value = "  +3.80%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.386"
$ws.Range("D9").Style = $plainStyle
$ws.Range("E9").Value = "  +2.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0807"
$ws.Range("D10").Style = $plainStyle
$ws.Range("E10").Value = "  +2.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.103"
$ws.Range("D11").Style = $plainStyle
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.11"
$ws.Range("D12").Style = $plainStyle
$ws.Range("E12").Value = "  +7.87%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.85"
$ws.Range("D13").Style = $plainStyle
$ws.Range("E13").Value = "  +7.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.852"
$ws.Range("D14").Style = $plainStyle
$ws.Range("E14").Value = "  +1.73%  "
$ws.Range("D15").Value = "2.297.21"
$ws.Range("E15").Value = "  +2.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.48"
$ws.Range("D16").Style = $plainStyle
$ws.Range("E16").Value = "  +3.74%  "
$ws.Range("D17").Value = "2.012.94"
$ws.Range("E17").Value = "  +3.11%  "
$ws.Range("D18").Value = "37.171.57"
$ws.Range("E18").Value = "  +1.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.45"
$ws.Range("D19").Style = $plainStyle
$ws.Range("E19").Value = "  +1.10%  "
$ws.Range("D20").Value = "0.0₃0865"
$ws.Range("E20").Value = "  +2.27%  "
$ws.Range("E21").Value = "  +3.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.93"
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("E24").Value = "  +1.59%  "
$ws.Range("E25").Value = "  +0.65%  "
$ws.Range("E26").Value = "  +3.74%  "
$ws.Range("E27").Value = "  +2.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "164.20"
$ws.Range("D28").Style = $plainStyle
$ws.Range("E28").Value = "  +2.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.71"
$ws.Range("D29").Style = $plainStyle
$ws.Range("E29").Value = "  +2.00%  "
$ws.Range("E30").Value = "  +13.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.122"
$ws.Range("D31").Style = $plainStyle
$ws.Range("E31").Value = "  +1.51%  "
$ws.Range("E32").Value = "  +2.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0661"
$ws.Range("D33").Style = $plainStyle
$ws.Range("E33").Value = "  +8.39%  "
$ws.Range("E34").Value = "  +3.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.40"
$ws.Range("D35").Style = $plainStyle
$ws.Range("E35").Value = "  +5.36%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  +2.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.30"
$ws.Range("D38").Style = $plainStyle
$ws.Range("E38").Value = "  -4.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.45"
$ws.Range("D39").Style = $plainStyle
$ws.Range("E39").Value = "  +2.37%  "
$ws.Range("E40").Value = "  +0.45%  "
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("E42").Value = "  +2.74%  "
$ws.Range("E43").Value = "  +2.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.84"
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = "  +7.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.34"
$ws.Range("D45").Style = $plainStyle
$ws.Range("E45").Value = "  +4.01%  "
$ws.Range("D46").Value = "1.380.80"
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("E47").Value = "  +2.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.39"
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = "  +3.79%  "
$ws.Range("E49").Value = "  +14.99%  "
$ws.Range("E50").Value = "  +0.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.53"
$ws.Range("D51").Style = $plainStyle
$ws.Range("E51").Value = "  +5.97%  "
